$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.310.72"
$ws.Range("E2").Value = "'  -0.99%  "
$ws.Range("D3").Value = "'2.594.66"
$ws.Range("E3").Value = "'  -2.94%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'509.37"
$ws.Range("E5").Value = "'  -0.42%  "
$ws.Range("D6").Value = "'153.92"
$ws.Range("E6").Value = "'  -2.11%  "
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("D9").Value = "'2.603.74"
$ws.Range("E9").Value = "'  -2.50%  "
$ws.Range("D10").Value = "'6.66"
$ws.Range("E10").Value = "'  +4.58%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "'  -1.12%  "
$ws.Range("E12").Value = "'  -0.91%  "
$ws.Range("E13").Value = "'  +1.71%  "
$ws.Range("D14").Value = "'3.048.99"
$ws.Range("E14").Value = "'  -2.82%  "
$ws.Range("D15").Value = "'60.320.54"
$ws.Range("E15").Value = "'  -1.08%  "
$ws.Range("D16").Value = "'21.49"
$ws.Range("E16").Value = "'  -1.39%  "
$ws.Range("E17").Value = "'  +0.08%  "
$ws.Range("D18").Value = "'2.595.07"
$ws.Range("E18").Value = "'  -2.67%  "
$ws.Range("D19").Value = "'4.74"
$ws.Range("E19").Value = "'  -1.34%  "
$ws.Range("D20").Value = "'352.74"
$ws.Range("E20").Value = "'  +1.30%  "
$ws.Range("D21").Value = "'10.53"
$ws.Range("E21").Value = "'  +0.22%  "
$ws.Range("E22").Value = "'  -0.82%  "
$ws.Range("E23").Value = "'  -0.05%  "
$ws.Range("D24").Value = "'60.35"
$ws.Range("E24").Value = "'  +0.19%  "
$ws.Range("D25").Value = "'0.420"
$ws.Range("E25").Value = "'  -0.54%  "
$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "'  -0.61%  "
$ws.Range("D28").Value = "'0.0₃0837"
$ws.Range("E28").Value = "'  -2.72%  "
$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = "'  -2.46%  "
$ws.Range("E30").Value = "'  +0.06%  "
$ws.Range("D31").Value = "'19.34"
$ws.Range("E31").Value = "'  -0.84%  "
$ws.Range("D32").Value = "'150.88"
$ws.Range("E32").Value = "'  -4.27%  "
$ws.Range("E33").Value = "'  -0.88%  "
$ws.Range("D34").Value = "'5.72"
$ws.Range("E34").Value = "'  +0.32%  "
$ws.Range("E35").Value = "'  -1.46%  "
$ws.Range("E36").Value = "'  -3.01%  "
$ws.Range("D37").Value = "'0.874"
$ws.Range("E37").Value = "'  +4.47%  "
$ws.Range("E38").Value = "'  -3.14%  "
$ws.Range("B39").Value = "'Fetch.AI"
$ws.Range("C39").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.840"
$ws.Range("E39").Value = "'  -2.30%  "
$ws.Range("B40").Value = "'OKB"
$ws.Range("C40").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'36.10"
$ws.Range("E40").Value = "'  +1.97%  "
$ws.Range("D42").Value = "'294.54"
$ws.Range("E42").Value = "'  -4.97%  "
$ws.Range("E43").Value = "'  -0.53%  "
$ws.Range("D44").Value = "'0.618"
$ws.Range("E44").Value = "'  -3.84%  "
$ws.Range("E45").Value = "'  -0.26%  "
$ws.Range("D46").Value = "'0.0554"
$ws.Range("E46").Value = "'  -4.20%  "
$ws.Range("D47").Value = "'19.60"
$ws.Range("E47").Value = "'  -1.44%  "
$ws.Range("B48").Value = "'RenderToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.86"
$ws.Range("E48").Value = "'  +0.40%  "
$ws.Range("B49").Value = "'VeChain"
$ws.Range("C49").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0234"
$ws.Range("E49").Value = "'  -1.14%  "
$ws.Range("E50").Value = "'  -0.12%  "
$ws.Range("D51").Value = "'1.987.94"
$ws.Range("E51").Value = "'  -2.71%  "
